$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('J6').Value = 55
$ws.Range('J7').Value = 122

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J6').Value = 74
$ws.Range('J7').Value = 204

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 239
$ws.Range('J3').Value = 232
$ws.Range('J6').Value = 244
$ws.Range('J7').Value = 766

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J6').Value = 603
$ws.Range('J7').Value = 1675

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J4').Value = 6
$ws.Range('J7').Value = 198

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('J6').Value = 59
$ws.Range('J7').Value = 144

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J6').Value = 146
$ws.Range('J7').Value = 326

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J4').Value = 122
$ws.Range('J6').Value = 204
$ws.Range('J7').Value = 766
$ws.Range('J8').Value = 1675
$ws.Range('J10').Value = 198
$ws.Range('J14').Value = 144
$ws.Range('J15').Value = 326
$ws.Range('J19').Value = 773
$ws.Range('J25').Value = 131
$ws.Range('J29').Value = 1432
$ws.Range('J31').Value = 273
$ws.Range('I33').Value = 1147
$ws.Range('J33').Value = 1212
$ws.Range('J34').Value = 124
$ws.Range('J36').Value = 362
$ws.Range('J37').Value = 826
$ws.Range('J41').Value = 196
$ws.Range('J42').Value = 1155
$ws.Range('J44').Value = 207
$ws.Range('J47').Value = 195
$ws.Range('J48').Value = 303
$ws.Range('J49').Value = 166
$ws.Range('J50').Value = 160
$ws.Range('J51').Value = 330
$ws.Range('J52').Value = 682
$ws.Range('J53').Value = 394
$ws.Range('J54').Value = 523
$ws.Range('J55').Value = 421
$ws.Range('J56').Value = 38
$ws.Range('J57').Value = 127
$ws.Range('J63').Value = 86
$ws.Range('J65').Value = 672
$ws.Range('J66').Value = 81
$ws.Range('G67').Value = 1203
$ws.Range('J67').Value = 1000
$ws.Range('J73').Value = 258
$ws.Range('J74').Value = 30
$ws.Range('J77').Value = 185
$ws.Range('J78').Value = 311
$ws.Range('J79').Value = 741
$ws.Range('J83').Value = 536
$ws.Range('J85').Value = 1105
$ws.Range('J90').Value = 282
$ws.Range('J91').Value = 309
$ws.Range('J92').Value = 85
$ws.Range('J93').Value = 113
$ws.Range('J96').Value = 291
$ws.Range('J97').Value = 244
$ws.Range('J98').Value = 199
$ws.Range('G101').Value = 24699
$ws.Range('I101').Value = 26233
$ws.Range('J101').Value = 26752

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J6').Value = 299
$ws.Range('J7').Value = 773

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 7047
$ws.Range('J3').Value = 7439
$ws.Range('G4').Value = 1474
$ws.Range('I4').Value = 1776
$ws.Range('J4').Value = 1625
$ws.Range('J5').Value = 583
$ws.Range('J6').Value = 10058
$ws.Range('G7').Value = 24699
$ws.Range('I7').Value = 26233
$ws.Range('J7').Value = 26752

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('J4').Value = 10
$ws.Range('J7').Value = 131

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 433
$ws.Range('J3').Value = 505
$ws.Range('J7').Value = 1432

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J6').Value = 96
$ws.Range('J7').Value = 273

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J3').Value = 399
$ws.Range('I4').Value = 52
$ws.Range('J6').Value = 430
$ws.Range('I7').Value = 1147
$ws.Range('J7').Value = 1212

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('J6').Value = 50
$ws.Range('J7').Value = 124

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J2').Value = 117
$ws.Range('J7').Value = 362

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 247
$ws.Range('J7').Value = 826

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J3').Value = 29
$ws.Range('J6').Value = 119
$ws.Range('J7').Value = 196

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J2').Value = 240
$ws.Range('J3').Value = 231
$ws.Range('J6').Value = 615
$ws.Range('J7').Value = 1155

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J3').Value = 46
$ws.Range('J4').Value = 13
$ws.Range('J6').Value = 82
$ws.Range('J7').Value = 207

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('J3').Value = 52
$ws.Range('J7').Value = 195

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J3').Value = 57
$ws.Range('J7').Value = 303

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('J3').Value = 31
$ws.Range('J6').Value = 93
$ws.Range('J7').Value = 166

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J3').Value = 40
$ws.Range('J7').Value = 160

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J3').Value = 88
$ws.Range('J7').Value = 330

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J3').Value = 195
$ws.Range('J4').Value = 25
$ws.Range('J7').Value = 682

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J6').Value = 263
$ws.Range('J7').Value = 394

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 130
$ws.Range('J4').Value = 41
$ws.Range('J6').Value = 243
$ws.Range('J7').Value = 523

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J6').Value = 238
$ws.Range('J7').Value = 421

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range('J6').Value = 18
$ws.Range('J7').Value = 38

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('J2').Value = 31
$ws.Range('J6').Value = 59
$ws.Range('J7').Value = 127

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 192
$ws.Range('J7').Value = 672

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('J6').Value = 46
$ws.Range('J7').Value = 81

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 255
$ws.Range('J3').Value = 374
$ws.Range('G4').Value = 50
$ws.Range('J6').Value = 277
$ws.Range('G7').Value = 1203
$ws.Range('J7').Value = 1000

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J2').Value = 83
$ws.Range('J6').Value = 95
$ws.Range('J7').Value = 258

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range('J3').Value = 10
$ws.Range('J7').Value = 30

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J4').Value = 18
$ws.Range('J7').Value = 185

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J2').Value = 83
$ws.Range('J7').Value = 311

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J6').Value = 221
$ws.Range('J7').Value = 741

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J2').Value = 156
$ws.Range('J7').Value = 536

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 297
$ws.Range('J7').Value = 1105

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J2').Value = 100
$ws.Range('J3').Value = 77
$ws.Range('J6').Value = 84
$ws.Range('J7').Value = 282

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J2').Value = 81
$ws.Range('J7').Value = 309

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('J6').Value = 27
$ws.Range('J7').Value = 85

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('J3').Value = 35
$ws.Range('J7').Value = 113

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J2').Value = 56
$ws.Range('J6').Value = 152

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J3').Value = 77
$ws.Range('J7').Value = 291

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J3').Value = 23
$ws.Range('J7').Value = 244

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J6').Value = 129
$ws.Range('J7').Value = 199
